# "import lai ds hoc sinh" -> append one more student id to the
# download_grade sheet (row 6, column A), same shape as the existing
# rows: a numeric-looking id stored as text (so Excel's "number stored as
# text" checker flags/ignores it just like rows 2-5 already do).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A6")
$cell.NumberFormat = "@"      # text format first, so the value below is kept as text, not coerced to a number
$cell.Value = "20000000"

# Worksheet's used range (and therefore its <dimension>) now naturally
# extends to A1:B6 as a side effect of populating A6.

# Best-effort: tell Excel's error checker to ignore the "number stored as
# text" warning for the whole updated range, mirroring the existing
# ignoredErrors/numberStoredAsText flag on A1:B5.
try {
    $errs = $ws.Range("A1:B6").Errors
    $errs.Item(9).Ignore = $true
} catch {
}
